# chore: adapt column header formatting to respective input file names (#7)
#
# Renames the header row suffixes ("_old" -> "_FV2410", "_new" -> "_FV2504"),
# wraps the data range in an Excel Table (Table1), and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row labels (row 1) ------------------------------------
$usedRange = $ws.UsedRange
$lastCol = $usedRange.Columns.Count
$lastRow = $usedRange.Rows.Count

for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -eq $null) { continue }
    $newVal = $val -replace '_old$', '_FV2410'
    $newVal = $newVal -replace '_new$', '_FV2504'
    if ($newVal -ne $val) {
        $cell.Value2 = $newVal
    }
}

# --- 2. Turn the data range into an Excel Table (Table1) --------------------
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3. Freeze the header row (row 1) ----------------------------------------
$ws.Range("A2").Select() | Out-Null
$null = ($ws.Application.ActiveWindow.FreezePanes = $true)
